# Inserts a new price-record row at row 77 (pushing the existing rows
# 77-138 down to 78-139, and growing the sheet from A1:R138 to A1:R139).
# The new row carries the same per-row formatting (date style on column D)
# as the rest of the table because Insert() shifts the formatted cells down
# and fills the freshly inserted row with the style of the row it displaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("77").Insert()

$ws.Range("A77").Value = 8
$ws.Range("B77").Value = "Terminal La Palmera de La Serena"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value2 = 44447
$ws.Range("E77").Value = 4
$ws.Range("F77").Value = 100112003
$ws.Range("G77").Value = "Ajo"
$ws.Range("H77").Value = "Chino"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 600
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = 15500
$ws.Range("N77").Value = "`$/caja 10 kilos"
$ws.Range("O77").Value = "China"
$ws.Range("P77").Value = 1550
$ws.Range("Q77").Value = 10
$ws.Range("R77").Value = "Hortaliza"
